$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'09/08/2023"
$ws.Range("B2").Value = "'4004.00"
$ws.Range("C2").Value = "'4004.00"
$ws.Range("D2").Value = "'5000.00"
$ws.Range("E2").Value = "'5000.00"
$ws.Range("F2").Value = "'996.00"
$ws.Range("G2").Value = "'124.88"

$ws.Range("A2:G2").ClearFormats()
